# version 4: add drop down menu of sub file name
#
# The mock-data table lists ID numbers (column A) and names (column B).
# Two sample records are removed from the table so that the remaining
# rows shift up (row count goes from 13 to 11):
#   - A673126970 / Name_17
#   - A726610453 / Name_10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$idsToRemove = @("A673126970", "A726610453")

foreach ($id in $idsToRemove) {
    $cell = $ws.Columns.Item(1).Find($id)
    if ($cell -ne $null) {
        $ws.Rows.Item($cell.Row).Delete()
    }
}
